# Auto-generated: adds rows 9-21 of species observation data to the Artfynd sheet
$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Row 9
$ws.Cells.Item(9, 1).Value = 131187727
$ws.Cells.Item(9, 2).Value = 57073
$ws.Cells.Item(9, 4).Value = "LC"
$ws.Cells.Item(9, 5).Value = 100138
$ws.Cells.Item(9, 6).Value = "Tjäder"
$ws.Cells.Item(9, 7).Value = "Tetrao urogallus"
$ws.Cells.Item(9, 8).Value = "Linnaeus, 1758"
$ws.Cells.Item(9, 9).Value = ""
$ws.Cells.Item(9, 11).Value = ""
$ws.Cells.Item(9, 12).Value = ""
$ws.Cells.Item(9, 13).Value = "färsk spillning"
$ws.Cells.Item(9, 14).Value = ""
$ws.Cells.Item(9, 16).Value = "Svartå, Dlr"
$ws.Cells.Item(9, 17).Value = 511543
$ws.Cells.Item(9, 18).Value = 6697760
$ws.Cells.Item(9, 19).Value = 25
$ws.Cells.Item(9, 20).Value = "Dalarna"
$ws.Cells.Item(9, 21).Value = "Borlänge"
$ws.Cells.Item(9, 22).Value = "Dalarna"
$ws.Cells.Item(9, 23).Value = "Stora Tuna"
$ws.Cells.Item(9, 25).Value = "'2026-02-16"
$ws.Cells.Item(9, 27).Value = "'2026-02-16"
$ws.Cells.Item(9, 30).Value = $false
$ws.Cells.Item(9, 31).Value = $false
$ws.Cells.Item(9, 33).Value = $false
$ws.Cells.Item(9, 46).Value = ""
$ws.Cells.Item(9, 49).Value = "Anna-Lena Thommson"
$ws.Cells.Item(9, 50).Value = "Anna-Lena Thommson"
$ws.Cells.Item(9, 51).Value = ""

# Row 10
$ws.Cells.Item(10, 1).Value = 131187799
$ws.Cells.Item(10, 2).Value = 79243
$ws.Cells.Item(10, 4).Value = "NT"
$ws.Cells.Item(10, 5).Value = 6425
$ws.Cells.Item(10, 6).Value = "Garnlav"
$ws.Cells.Item(10, 7).Value = "Alectoria sarmentosa"
$ws.Cells.Item(10, 8).Value = "(Ach.) Ach."
$ws.Cells.Item(10, 9).Value = ""
$ws.Cells.Item(10, 10).Value = ""
$ws.Cells.Item(10, 11).Value = ""
$ws.Cells.Item(10, 14).Value = ""
$ws.Cells.Item(10, 16).Value = "Svartå, Dlr"
$ws.Cells.Item(10, 17).Value = 511308
$ws.Cells.Item(10, 18).Value = 6697583
$ws.Cells.Item(10, 19).Value = 25
$ws.Cells.Item(10, 20).Value = "Dalarna"
$ws.Cells.Item(10, 21).Value = "Borlänge"
$ws.Cells.Item(10, 22).Value = "Dalarna"
$ws.Cells.Item(10, 23).Value = "Stora Tuna"
$ws.Cells.Item(10, 25).Value = "'2026-02-16"
$ws.Cells.Item(10, 27).Value = "'2026-02-16"
$ws.Cells.Item(10, 29).Value = "På gran."
$ws.Cells.Item(10, 30).Value = $false
$ws.Cells.Item(10, 31).Value = $false
$ws.Cells.Item(10, 32).Value = ""
$ws.Cells.Item(10, 33).Value = $false
$ws.Cells.Item(10, 46).Value = ""
$ws.Cells.Item(10, 49).Value = "Anna-Lena Thommson"
$ws.Cells.Item(10, 50).Value = "Anna-Lena Thommson"
$ws.Cells.Item(10, 51).Value = ""

# Row 11
$ws.Cells.Item(11, 1).Value = 131191997
$ws.Cells.Item(11, 2).Value = 57884
$ws.Cells.Item(11, 4).Value = "NT"
$ws.Cells.Item(11, 5).Value = 100109
$ws.Cells.Item(11, 6).Value = "Tretåig hackspett"
$ws.Cells.Item(11, 7).Value = "Picoides tridactylus"
$ws.Cells.Item(11, 8).Value = "(Linnaeus, 1758)"
$ws.Cells.Item(11, 9).Value = ""
$ws.Cells.Item(11, 11).Value = ""
$ws.Cells.Item(11, 12).Value = ""
$ws.Cells.Item(11, 13).Value = "äldre spår"
$ws.Cells.Item(11, 14).Value = ""
$ws.Cells.Item(11, 16).Value = "Svartå, Dlr"
$ws.Cells.Item(11, 17).Value = 511355
$ws.Cells.Item(11, 18).Value = 6697418
$ws.Cells.Item(11, 19).Value = 10
$ws.Cells.Item(11, 20).Value = "Dalarna"
$ws.Cells.Item(11, 21).Value = "Borlänge"
$ws.Cells.Item(11, 22).Value = "Dalarna"
$ws.Cells.Item(11, 23).Value = "Stora Tuna"
$ws.Cells.Item(11, 25).Value = "'2026-02-16"
$ws.Cells.Item(11, 26).Value = "12:09"
$ws.Cells.Item(11, 27).Value = "'2026-02-16"
$ws.Cells.Item(11, 28).Value = "12:09"
$ws.Cells.Item(11, 29).Value = "Riinghack äldre på tall."
$ws.Cells.Item(11, 30).Value = $false
$ws.Cells.Item(11, 31).Value = $false
$ws.Cells.Item(11, 33).Value = $false
$ws.Cells.Item(11, 46).Value = ""
$ws.Cells.Item(11, 49).Value = "Lars-Erik Nilsson"
$ws.Cells.Item(11, 50).Value = "Lars-Erik Nilsson, Anna-Lena Thommson"
$ws.Cells.Item(11, 51).Value = ""

# Row 12
$ws.Cells.Item(12, 1).Value = 131187835
$ws.Cells.Item(12, 2).Value = 57073
$ws.Cells.Item(12, 4).Value = "LC"
$ws.Cells.Item(12, 5).Value = 100138
$ws.Cells.Item(12, 6).Value = "Tjäder"
$ws.Cells.Item(12, 7).Value = "Tetrao urogallus"
$ws.Cells.Item(12, 8).Value = "Linnaeus, 1758"
$ws.Cells.Item(12, 9).Value = ""
$ws.Cells.Item(12, 11).Value = ""
$ws.Cells.Item(12, 12).Value = ""
$ws.Cells.Item(12, 13).Value = "äldre spår"
$ws.Cells.Item(12, 14).Value = ""
$ws.Cells.Item(12, 16).Value = "Svartå, Dlr"
$ws.Cells.Item(12, 17).Value = 511382
$ws.Cells.Item(12, 18).Value = 6697458
$ws.Cells.Item(12, 19).Value = 25
$ws.Cells.Item(12, 20).Value = "Dalarna"
$ws.Cells.Item(12, 21).Value = "Borlänge"
$ws.Cells.Item(12, 22).Value = "Dalarna"
$ws.Cells.Item(12, 23).Value = "Stora Tuna"
$ws.Cells.Item(12, 25).Value = "'2026-02-16"
$ws.Cells.Item(12, 27).Value = "'2026-02-16"
$ws.Cells.Item(12, 29).Value = "Betad tallkrona."
$ws.Cells.Item(12, 30).Value = $false
$ws.Cells.Item(12, 31).Value = $false
$ws.Cells.Item(12, 33).Value = $false
$ws.Cells.Item(12, 46).Value = ""
$ws.Cells.Item(12, 49).Value = "Anna-Lena Thommson"
$ws.Cells.Item(12, 50).Value = "Anna-Lena Thommson"
$ws.Cells.Item(12, 51).Value = ""

# Row 13
$ws.Cells.Item(13, 1).Value = 131187780
$ws.Cells.Item(13, 2).Value = 79243
$ws.Cells.Item(13, 4).Value = "NT"
$ws.Cells.Item(13, 5).Value = 6425
$ws.Cells.Item(13, 6).Value = "Garnlav"
$ws.Cells.Item(13, 7).Value = "Alectoria sarmentosa"
$ws.Cells.Item(13, 8).Value = "(Ach.) Ach."
$ws.Cells.Item(13, 9).Value = ""
$ws.Cells.Item(13, 10).Value = ""
$ws.Cells.Item(13, 11).Value = ""
$ws.Cells.Item(13, 14).Value = ""
$ws.Cells.Item(13, 16).Value = "Svatå, Dlr"
$ws.Cells.Item(13, 17).Value = 511335
$ws.Cells.Item(13, 18).Value = 6697864
$ws.Cells.Item(13, 19).Value = 25
$ws.Cells.Item(13, 20).Value = "Dalarna"
$ws.Cells.Item(13, 21).Value = "Borlänge"
$ws.Cells.Item(13, 22).Value = "Dalarna"
$ws.Cells.Item(13, 23).Value = "Stora Tuna"
$ws.Cells.Item(13, 25).Value = "'2026-02-16"
$ws.Cells.Item(13, 27).Value = "'2026-02-16"
$ws.Cells.Item(13, 29).Value = "På gran."
$ws.Cells.Item(13, 30).Value = $false
$ws.Cells.Item(13, 31).Value = $false
$ws.Cells.Item(13, 32).Value = ""
$ws.Cells.Item(13, 33).Value = $false
$ws.Cells.Item(13, 46).Value = ""
$ws.Cells.Item(13, 49).Value = "Anna-Lena Thommson"
$ws.Cells.Item(13, 50).Value = "Anna-Lena Thommson"
$ws.Cells.Item(13, 51).Value = ""

# Row 14
$ws.Cells.Item(14, 1).Value = 131191949
$ws.Cells.Item(14, 2).Value = 79243
$ws.Cells.Item(14, 4).Value = "NT"
$ws.Cells.Item(14, 5).Value = 6425
$ws.Cells.Item(14, 6).Value = "Garnlav"
$ws.Cells.Item(14, 7).Value = "Alectoria sarmentosa"
$ws.Cells.Item(14, 8).Value = "(Ach.) Ach."
$ws.Cells.Item(14, 9).Value = ""
$ws.Cells.Item(14, 10).Value = ""
$ws.Cells.Item(14, 11).Value = ""
$ws.Cells.Item(14, 14).Value = ""
$ws.Cells.Item(14, 16).Value = "Svartå, Dlr"
$ws.Cells.Item(14, 17).Value = 511393
$ws.Cells.Item(14, 18).Value = 6697824
$ws.Cells.Item(14, 19).Value = 10
$ws.Cells.Item(14, 20).Value = "Dalarna"
$ws.Cells.Item(14, 21).Value = "Borlänge"
$ws.Cells.Item(14, 22).Value = "Dalarna"
$ws.Cells.Item(14, 23).Value = "Stora Tuna"
$ws.Cells.Item(14, 25).Value = "'2026-02-16"
$ws.Cells.Item(14, 26).Value = "10:33"
$ws.Cells.Item(14, 27).Value = "'2026-02-16"
$ws.Cells.Item(14, 28).Value = "10:33"
$ws.Cells.Item(14, 30).Value = $false
$ws.Cells.Item(14, 31).Value = $false
$ws.Cells.Item(14, 32).Value = ""
$ws.Cells.Item(14, 33).Value = $false
$ws.Cells.Item(14, 46).Value = ""
$ws.Cells.Item(14, 49).Value = "Lars-Erik Nilsson"
$ws.Cells.Item(14, 50).Value = "Lars-Erik Nilsson, Anna-Lena Thommson"
$ws.Cells.Item(14, 51).Value = ""

# Row 15
$ws.Cells.Item(15, 1).Value = 131192467
$ws.Cells.Item(15, 2).Value = 57073
$ws.Cells.Item(15, 4).Value = "LC"
$ws.Cells.Item(15, 5).Value = 100138
$ws.Cells.Item(15, 6).Value = "Tjäder"
$ws.Cells.Item(15, 7).Value = "Tetrao urogallus"
$ws.Cells.Item(15, 8).Value = "Linnaeus, 1758"
$ws.Cells.Item(15, 9).Value = ""
$ws.Cells.Item(15, 11).Value = ""
$ws.Cells.Item(15, 12).Value = ""
$ws.Cells.Item(15, 13).Value = ""
$ws.Cells.Item(15, 14).Value = ""
$ws.Cells.Item(15, 16).Value = "Svartå, Dlr"
$ws.Cells.Item(15, 17).Value = 511370
$ws.Cells.Item(15, 18).Value = 6697492
$ws.Cells.Item(15, 19).Value = 25
$ws.Cells.Item(15, 20).Value = "Dalarna"
$ws.Cells.Item(15, 21).Value = "Borlänge"
$ws.Cells.Item(15, 22).Value = "Dalarna"
$ws.Cells.Item(15, 23).Value = "Stora Tuna"
$ws.Cells.Item(15, 25).Value = "'2026-02-16"
$ws.Cells.Item(15, 26).Value = "11:38"
$ws.Cells.Item(15, 27).Value = "'2026-02-16"
$ws.Cells.Item(15, 28).Value = "11:38"
$ws.Cells.Item(15, 29).Value = "Två tjäderbetade tallar,"
$ws.Cells.Item(15, 30).Value = $false
$ws.Cells.Item(15, 31).Value = $false
$ws.Cells.Item(15, 33).Value = $false
$ws.Cells.Item(15, 46).Value = ""
$ws.Cells.Item(15, 49).Value = "Lars-Erik Nilsson"
$ws.Cells.Item(15, 50).Value = "Lars-Erik Nilsson, Anna-Lena Thommson"
$ws.Cells.Item(15, 51).Value = ""

# Row 16
$ws.Cells.Item(16, 1).Value = 131187861
$ws.Cells.Item(16, 2).Value = 58043
$ws.Cells.Item(16, 4).Value = "NT"
$ws.Cells.Item(16, 5).Value = 103021
$ws.Cells.Item(16, 6).Value = "Talltita"
$ws.Cells.Item(16, 7).Value = "Poecile montanus"
$ws.Cells.Item(16, 8).Value = "(Conrad von Baldenstein, 1827)"
$ws.Cells.Item(16, 9).Value = "'1"
$ws.Cells.Item(16, 11).Value = ""
$ws.Cells.Item(16, 12).Value = ""
$ws.Cells.Item(16, 13).Value = "födosökande"
$ws.Cells.Item(16, 14).Value = ""
$ws.Cells.Item(16, 16).Value = "Svartå, Dlr"
$ws.Cells.Item(16, 17).Value = 511287
$ws.Cells.Item(16, 18).Value = 6697400
$ws.Cells.Item(16, 19).Value = 25
$ws.Cells.Item(16, 20).Value = "Dalarna"
$ws.Cells.Item(16, 21).Value = "Borlänge"
$ws.Cells.Item(16, 22).Value = "Dalarna"
$ws.Cells.Item(16, 23).Value = "Stora Tuna"
$ws.Cells.Item(16, 25).Value = "'2026-02-16"
$ws.Cells.Item(16, 27).Value = "'2026-02-16"
$ws.Cells.Item(16, 30).Value = $false
$ws.Cells.Item(16, 31).Value = $false
$ws.Cells.Item(16, 33).Value = $false
$ws.Cells.Item(16, 46).Value = ""
$ws.Cells.Item(16, 49).Value = "Anna-Lena Thommson"
$ws.Cells.Item(16, 50).Value = "Anna-Lena Thommson"
$ws.Cells.Item(16, 51).Value = ""

# Row 17
$ws.Cells.Item(17, 1).Value = 131187741
$ws.Cells.Item(17, 2).Value = 79243
$ws.Cells.Item(17, 4).Value = "NT"
$ws.Cells.Item(17, 5).Value = 6425
$ws.Cells.Item(17, 6).Value = "Garnlav"
$ws.Cells.Item(17, 7).Value = "Alectoria sarmentosa"
$ws.Cells.Item(17, 8).Value = "(Ach.) Ach."
$ws.Cells.Item(17, 9).Value = ""
$ws.Cells.Item(17, 10).Value = ""
$ws.Cells.Item(17, 11).Value = ""
$ws.Cells.Item(17, 14).Value = ""
$ws.Cells.Item(17, 16).Value = "Svartå, Dlr"
$ws.Cells.Item(17, 17).Value = 511476
$ws.Cells.Item(17, 18).Value = 6697750
$ws.Cells.Item(17, 19).Value = 25
$ws.Cells.Item(17, 20).Value = "Dalarna"
$ws.Cells.Item(17, 21).Value = "Borlänge"
$ws.Cells.Item(17, 22).Value = "Dalarna"
$ws.Cells.Item(17, 23).Value = "Stora Tuna"
$ws.Cells.Item(17, 25).Value = "'2026-02-16"
$ws.Cells.Item(17, 27).Value = "'2026-02-16"
$ws.Cells.Item(17, 29).Value = "Rikligt på äldre tall."
$ws.Cells.Item(17, 30).Value = $false
$ws.Cells.Item(17, 31).Value = $false
$ws.Cells.Item(17, 32).Value = ""
$ws.Cells.Item(17, 33).Value = $false
$ws.Cells.Item(17, 46).Value = ""
$ws.Cells.Item(17, 49).Value = "Anna-Lena Thommson"
$ws.Cells.Item(17, 50).Value = "Anna-Lena Thommson"
$ws.Cells.Item(17, 51).Value = ""

# Row 18
$ws.Cells.Item(18, 1).Value = 131187791
$ws.Cells.Item(18, 2).Value = 57073
$ws.Cells.Item(18, 4).Value = "LC"
$ws.Cells.Item(18, 5).Value = 100138
$ws.Cells.Item(18, 6).Value = "Tjäder"
$ws.Cells.Item(18, 7).Value = "Tetrao urogallus"
$ws.Cells.Item(18, 8).Value = "Linnaeus, 1758"
$ws.Cells.Item(18, 9).Value = ""
$ws.Cells.Item(18, 11).Value = ""
$ws.Cells.Item(18, 12).Value = ""
$ws.Cells.Item(18, 13).Value = "färsk spillning"
$ws.Cells.Item(18, 14).Value = ""
$ws.Cells.Item(18, 16).Value = "Svatå, Dlr"
$ws.Cells.Item(18, 17).Value = 511301
$ws.Cells.Item(18, 18).Value = 6697864
$ws.Cells.Item(18, 19).Value = 25
$ws.Cells.Item(18, 20).Value = "Dalarna"
$ws.Cells.Item(18, 21).Value = "Borlänge"
$ws.Cells.Item(18, 22).Value = "Dalarna"
$ws.Cells.Item(18, 23).Value = "Stora Tuna"
$ws.Cells.Item(18, 25).Value = "'2026-02-16"
$ws.Cells.Item(18, 27).Value = "'2026-02-16"
$ws.Cells.Item(18, 30).Value = $false
$ws.Cells.Item(18, 31).Value = $false
$ws.Cells.Item(18, 33).Value = $false
$ws.Cells.Item(18, 46).Value = ""
$ws.Cells.Item(18, 49).Value = "Anna-Lena Thommson"
$ws.Cells.Item(18, 50).Value = "Anna-Lena Thommson"
$ws.Cells.Item(18, 51).Value = ""

# Row 19
$ws.Cells.Item(19, 1).Value = 131187762
$ws.Cells.Item(19, 2).Value = 79243
$ws.Cells.Item(19, 4).Value = "NT"
$ws.Cells.Item(19, 5).Value = 6425
$ws.Cells.Item(19, 6).Value = "Garnlav"
$ws.Cells.Item(19, 7).Value = "Alectoria sarmentosa"
$ws.Cells.Item(19, 8).Value = "(Ach.) Ach."
$ws.Cells.Item(19, 9).Value = ""
$ws.Cells.Item(19, 10).Value = ""
$ws.Cells.Item(19, 11).Value = ""
$ws.Cells.Item(19, 14).Value = ""
$ws.Cells.Item(19, 16).Value = "Svartå, Dlr"
$ws.Cells.Item(19, 17).Value = 511511
$ws.Cells.Item(19, 18).Value = 6697866
$ws.Cells.Item(19, 19).Value = 25
$ws.Cells.Item(19, 20).Value = "Dalarna"
$ws.Cells.Item(19, 21).Value = "Borlänge"
$ws.Cells.Item(19, 22).Value = "Dalarna"
$ws.Cells.Item(19, 23).Value = "Stora Tuna"
$ws.Cells.Item(19, 25).Value = "'2026-02-16"
$ws.Cells.Item(19, 27).Value = "'2026-02-16"
$ws.Cells.Item(19, 29).Value = "På äldre tall."
$ws.Cells.Item(19, 30).Value = $false
$ws.Cells.Item(19, 31).Value = $false
$ws.Cells.Item(19, 32).Value = ""
$ws.Cells.Item(19, 33).Value = $false
$ws.Cells.Item(19, 46).Value = ""
$ws.Cells.Item(19, 49).Value = "Anna-Lena Thommson"
$ws.Cells.Item(19, 50).Value = "Anna-Lena Thommson"
$ws.Cells.Item(19, 51).Value = ""

# Row 20
$ws.Cells.Item(20, 1).Value = 131191884
$ws.Cells.Item(20, 2).Value = 79243
$ws.Cells.Item(20, 4).Value = "NT"
$ws.Cells.Item(20, 5).Value = 6425
$ws.Cells.Item(20, 6).Value = "Garnlav"
$ws.Cells.Item(20, 7).Value = "Alectoria sarmentosa"
$ws.Cells.Item(20, 8).Value = "(Ach.) Ach."
$ws.Cells.Item(20, 9).Value = ""
$ws.Cells.Item(20, 10).Value = ""
$ws.Cells.Item(20, 11).Value = ""
$ws.Cells.Item(20, 14).Value = ""
$ws.Cells.Item(20, 16).Value = "Svartå, Dlr"
$ws.Cells.Item(20, 17).Value = 511360
$ws.Cells.Item(20, 18).Value = 6697921
$ws.Cells.Item(20, 19).Value = 10
$ws.Cells.Item(20, 20).Value = "Dalarna"
$ws.Cells.Item(20, 21).Value = "Borlänge"
$ws.Cells.Item(20, 22).Value = "Dalarna"
$ws.Cells.Item(20, 23).Value = "Stora Tuna"
$ws.Cells.Item(20, 25).Value = "'2026-02-16"
$ws.Cells.Item(20, 26).Value = "10:24"
$ws.Cells.Item(20, 27).Value = "'2026-02-16"
$ws.Cells.Item(20, 28).Value = "10:24"
$ws.Cells.Item(20, 30).Value = $false
$ws.Cells.Item(20, 31).Value = $false
$ws.Cells.Item(20, 32).Value = ""
$ws.Cells.Item(20, 33).Value = $false
$ws.Cells.Item(20, 46).Value = ""
$ws.Cells.Item(20, 49).Value = "Lars-Erik Nilsson"
$ws.Cells.Item(20, 50).Value = "Lars-Erik Nilsson, Anna-Lena Thommson"
$ws.Cells.Item(20, 51).Value = ""

# Row 21
$ws.Cells.Item(21, 1).Value = 131191374
$ws.Cells.Item(21, 2).Value = 57884
$ws.Cells.Item(21, 4).Value = "NT"
$ws.Cells.Item(21, 5).Value = 100109
$ws.Cells.Item(21, 6).Value = "Tretåig hackspett"
$ws.Cells.Item(21, 7).Value = "Picoides tridactylus"
$ws.Cells.Item(21, 8).Value = "(Linnaeus, 1758)"
$ws.Cells.Item(21, 9).Value = ""
$ws.Cells.Item(21, 11).Value = ""
$ws.Cells.Item(21, 12).Value = ""
$ws.Cells.Item(21, 13).Value = "äldre spår"
$ws.Cells.Item(21, 14).Value = ""
$ws.Cells.Item(21, 16).Value = "Svartå, Dlr"
$ws.Cells.Item(21, 17).Value = 511332
$ws.Cells.Item(21, 18).Value = 6697755
$ws.Cells.Item(21, 19).Value = 10
$ws.Cells.Item(21, 20).Value = "Dalarna"
$ws.Cells.Item(21, 21).Value = "Borlänge"
$ws.Cells.Item(21, 22).Value = "Dalarna"
$ws.Cells.Item(21, 23).Value = "Stora Tuna"
$ws.Cells.Item(21, 25).Value = "'2026-02-16"
$ws.Cells.Item(21, 26).Value = "10:51"
$ws.Cells.Item(21, 27).Value = "'2026-02-16"
$ws.Cells.Item(21, 28).Value = "10:51"
$ws.Cells.Item(21, 29).Value = "Ringhack."
$ws.Cells.Item(21, 30).Value = $false
$ws.Cells.Item(21, 31).Value = $false
$ws.Cells.Item(21, 33).Value = $false
$ws.Cells.Item(21, 46).Value = ""
$ws.Cells.Item(21, 49).Value = "Lars-Erik Nilsson"
$ws.Cells.Item(21, 50).Value = "Lars-Erik Nilsson, Anna-Lena Thommson"
$ws.Cells.Item(21, 51).Value = ""

Write-Host "Added rows 9-21 with species observation data"